$wb = $excel.ActiveWorkbook

# --- Sheet "Jon": shift all years in A2:A36 up by one (2020->2021 ... 2054->2055),
#     then drop the old last row (r=37, year 2055) which is now a duplicate. ---
$wsJon = $wb.Worksheets.Item("Jon")
for ($r = 2; $r -le 36; $r++) {
    $yr = $wsJon.Cells.Item($r, 1).Value2
    $wsJon.Cells.Item($r, 1).Value2 = $yr + 1
}
$wsJon.Rows.Item(37).Delete() | Out-Null
$wsJon.Range("A2:A36").Select() | Out-Null

# --- Sheet "Jane": shift all years in A2:A39 up by one (2020->2021 ... 2057->2058),
#     then drop the old last row (r=40, year 2058) which is now a duplicate. ---
$wsJane = $wb.Worksheets.Item("Jane")
for ($r = 2; $r -le 39; $r++) {
    $yr = $wsJane.Cells.Item($r, 1).Value2
    $wsJane.Cells.Item($r, 1).Value2 = $yr + 1
}
$wsJane.Rows.Item(40).Delete() | Out-Null
$wsJane.Range("B3").Select() | Out-Null

# Jane becomes the active sheet/tab (was Debts).
$wsJane.Activate() | Out-Null
